# xlOil_Utils.xlsx - add "array functionality in xloSplit" test block to the
# Concat-Split sheet (sheet2): three new sections -
#   - "Split array of string"            (B25:H30)
#   - "Split horizontal array of string"  (B32:G37)
#   - "Ignore non-strings"                (B39:E42)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Block-Fill-FillNA (donor for style s=2)
$ws2 = $wb.Worksheets.Item(2)   # Concat-Split       (the sheet being edited)

$xlPasteFormats = -4122

function Copy-Style($srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Section header: "Split array of string"
# ---------------------------------------------------------------------------
$ws2.Range("B25").Value = "Split array of string"

# ---- Row 27: array formula anchor + its expected (test) values -----------
$ws2.Range("B27").Value = "Foo"
$ws2.Range("C27").Value = "#N/A"
$ws2.Range("E27").Value = "Foo"
$ws2.Range("F27").Value = "#N/A"

# ---- Row 28 ----------------------------------------------------------------
$ws2.Range("B28").Value = "Gr"
$ws2.Range("C28").Value = "n"
$ws2.Range("E28").Value = "Gr"
$ws2.Range("F28").Value = "n"

# ---- Row 29 ----------------------------------------------------------------
$ws2.Range("B29").Value = "Eggs"
$ws2.Range("C29").Value = "#N/A"
$ws2.Range("E29").Value = "Eggs"
$ws2.Range("F29").Value = "#N/A"

# ---- Row 30 ----------------------------------------------------------------
$ws2.Range("B30").Value = "Ham"
$ws2.Range("C30").Value = "#N/A"
$ws2.Range("E30").Value = "Ham"
$ws2.Range("F30").Value = "#N/A"

# Apply the "split-result" fill style (s=5) to B27:C30, then enter the single
# legacy-CSE array formula that spans the whole block (matches ref=B27:C30)
Copy-Style $ws2.Range("B5") $ws2.Range("B27:C30")
$ws2.Range("B27:C30").FormulaArray = '=_xll.xloSplit(I5:I8,"e")'

# Comparison / ISNA check columns (G:H), style s=16
Copy-Style $ws2.Range("B23") $ws2.Range("G27:H30")
$ws2.Range("G27").Formula = "=E27=B27"
$ws2.Range("H27").Formula = "=ISNA(C27)"
$ws2.Range("G28:G30").FormulaR1C1 = "=RC[-2]=RC[-5]"
$ws2.Range("H28").Formula = "=F28=C28"
$ws2.Range("H29").Formula = "=ISNA(C29)"
$ws2.Range("H30").Formula = "=ISNA(C30)"

# ---------------------------------------------------------------------------
# Section header: "Split horizontal array of string"
# ---------------------------------------------------------------------------
$ws2.Range("B32").Value = "Split horizontal array of string"

# ---- Row 34 / 35 - expected values -----------------------------------------
$ws2.Range("B34").Value = "Foo"
$ws2.Range("C34").Value = "B"
$ws2.Range("D34").Value = "B"
$ws2.Range("B35").Value = "#N/A"
$ws2.Range("C35").Value = "r"
$ws2.Range("D35").Value = "z"

Copy-Style $ws2.Range("B5") $ws2.Range("B34:D35")
$ws2.Range("B34:D35").FormulaArray = '=_xll.xloSplit(I5:K5,"a")'

# ---- Row 36 / 37 - reference data being compared against ------------------
$ws2.Range("B36").Value = "Foo"
$ws2.Range("C36").Value = "B"
$ws2.Range("D36").Value = "B"
$ws2.Range("B37").Value = "#N/A"
$ws2.Range("C37").Value = "r"
$ws2.Range("D37").Value = "z"

# Comparison columns (E:G), style s=16
Copy-Style $ws2.Range("B23") $ws2.Range("E34:G35")
$ws2.Range("E34").Formula = "=B36=B34"
$ws2.Range("F34").Formula = "=C36=C34"
$ws2.Range("G34").Formula = "=D36=D34"
$ws2.Range("E35").Formula = "=ISNA(B37)"
$ws2.Range("F35").Formula = "=C37=C35"
$ws2.Range("G35").Formula = "=D37=D35"

# ---------------------------------------------------------------------------
# Section header: "Ignore non-strings"
# ---------------------------------------------------------------------------
$ws2.Range("B39").Value = "Ignore non-strings"

# ---- Rows 40:42 - numeric input column (B) is left untouched by xloSplit --
Copy-Style $ws1.Range("C6") $ws2.Range("B40:B42")
$ws2.Range("B40").Value = 1
$ws2.Range("B41").Value = 2
$ws2.Range("B42").Value = 3

Copy-Style $ws2.Range("B5") $ws2.Range("C40:C42")
$ws2.Range("C41").Value = 2
$ws2.Range("C42").Value = 3
$ws2.Range("C40:C42").FormulaArray = '=_xll.xloSplit(B40:B42,"x")'

Copy-Style $ws2.Range("B23") $ws2.Range("E40:E42")
$ws2.Range("E40").Formula = "=C40=B40"
$ws2.Range("E41:E42").FormulaR1C1 = "=RC[-2]=RC[-3]"

# ---------------------------------------------------------------------------
# UI state: the user ends up with Concat-Split as the active sheet, selection
# sitting on G18
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("G18").Select() | Out-Null

Write-Host "xloSplit array-functionality test block added"
